$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: second date value, same style as B2
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C2").Value = 43986

# Row 5: new task "Generación de gráficas" using same style as A4 (task rows)
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "Generación de gráficas"

# Row 6: new section header "Modo entrenamiento" using same style as A3 (section header)
$ws.Range("A3").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "Modo entrenamiento"

# Row 7: new task "Entrenamiento IA vs IA básico" using same style as A4
$ws.Range("A4").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "Entrenamiento IA vs IA básico"

# C5: "3 h." using same style as B4 (hours column)
$ws.Range("B4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "3 h."

# C7: "1.5 h." using same style as B4
$ws.Range("B4").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = "1.5 h."

$excel.CutCopyMode = 0

# Update selection to match target state
$ws.Range("F6").Select()
